$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New bold header labels (row 18) - "Mean increase" in D18, "Median increase" in F18
$ws.Range("D18").Value = "Mean increase"
$ws.Range("D18").Font.Bold = $true

$ws.Range("F18").Value = "Median increase"
$ws.Range("F18").Font.Bold = $true

# New formulas (row 19) computing the percentage increase relative to a baseline
$ws.Range("D19").Formula = "=((E3 / 95.321842) * 100) - 100"
$ws.Range("D19").ClearFormats()

$ws.Range("F19").Formula = "=((E10 / 95.22216) * 100) - 100"
$ws.Range("F19").ClearFormats()

# Update the active selection to match the authored workbook state
[void]$ws.Range("D20").Select()
